# Daily attendance processing - 2026-01-22 14:02:06
# Applies the attendance-recording update described in the commit:
#  - Class Statistics totals (Recorded/Missing/Coverage %/Avg Attendance %)
#  - Per-group statistics (Recorded/Missing counts + Coverage %/Avg Attendance %)
#  - 24 session rows flip from "Not Recorded" (pink) to "Recorded" (green),
#    picking up a "System" recorder and full attendance (H = "total/total"
#    or a partial count where noted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write literal text that LOOKS like a percentage (e.g. "84.6%")
# without Excel's General-format auto-conversion turning it into a numeric
# percentage. We stage the text in a scratch cell formatted as Text, copy
# it, then PasteSpecial *values only* onto the destination so the
# destination keeps its own (existing) style/number format untouched.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

function Set-PctText($addr, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# ---------------------------------------------------------------------------
# Helper: flip a session detail row (columns A:I) from "Not Recorded" to
# "Recorded" - copies the green "Recorded" formatting from row 2 (A2:I2,
# a known Recorded row) onto the target row, then fills in Recorded By /
# Students / Status.
# ---------------------------------------------------------------------------
function Set-RecordedRow($row, $studentsText) {
    $ws.Range("A2:I2").Copy()
    $ws.Range("A" + $row + ":I" + $row).PasteSpecial(-4122)
    $ws.Range("G" + $row).Value = "System"
    $ws.Range("H" + $row).Value = $studentsText
    $ws.Range("I" + $row).Value = "Recorded"
}

# ---------------------------------------------------------------------------
# Class Statistics (K4:L10 block)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 273
$ws.Range("L7").Value = 45
Set-PctText "L9" "85.8%"
Set-PctText "L10" "77.7%"

# ---------------------------------------------------------------------------
# Group Statistics (rows 15-26): Recorded(O) / Missing(P) / Coverage%(R) /
# Avg Attendance%(S)
# ---------------------------------------------------------------------------
$ws.Range("O15").Value = 22
$ws.Range("P15").Value = 4
Set-PctText "R15" "84.6%"
Set-PctText "S15" "80.8%"

$ws.Range("O16").Value = 23
$ws.Range("P16").Value = 3
Set-PctText "R16" "88.5%"
Set-PctText "S16" "81.2%"

$ws.Range("O17").Value = 23
$ws.Range("P17").Value = 3
Set-PctText "R17" "88.5%"
Set-PctText "S17" "73.7%"

$ws.Range("O18").Value = 23
$ws.Range("P18").Value = 3
Set-PctText "R18" "88.5%"
Set-PctText "S18" "78.3%"

$ws.Range("O19").Value = 23
$ws.Range("P19").Value = 3
Set-PctText "R19" "88.5%"
Set-PctText "S19" "78.0%"

$ws.Range("O20").Value = 22
$ws.Range("P20").Value = 4
Set-PctText "R20" "84.6%"
Set-PctText "S20" "80.0%"

$ws.Range("O21").Value = 23
$ws.Range("P21").Value = 4
Set-PctText "R21" "85.2%"
Set-PctText "S21" "81.3%"

$ws.Range("O22").Value = 23
$ws.Range("P22").Value = 4
Set-PctText "R22" "85.2%"
Set-PctText "S22" "79.9%"

$ws.Range("O23").Value = 23
$ws.Range("P23").Value = 4
Set-PctText "R23" "85.2%"
Set-PctText "S23" "80.1%"

$ws.Range("O24").Value = 22
$ws.Range("P24").Value = 5
Set-PctText "R24" "81.5%"
Set-PctText "S24" "75.3%"

$ws.Range("O25").Value = 23
$ws.Range("P25").Value = 4
Set-PctText "R25" "85.2%"
Set-PctText "S25" "74.6%"

$ws.Range("O26").Value = 23
$ws.Range("P26").Value = 4
Set-PctText "R26" "85.2%"
Set-PctText "S26" "69.1%"

# ---------------------------------------------------------------------------
# Session detail rows flipping from "Not Recorded" -> "Recorded"
# ---------------------------------------------------------------------------
Set-RecordedRow 21 "26/26"
Set-RecordedRow 47 "27/27"
Set-RecordedRow 73 "26/26"
Set-RecordedRow 99 "27/27"
Set-RecordedRow 125 "30/30"
Set-RecordedRow 151 "23/23"
Set-RecordedRow 175 "22/23"
Set-RecordedRow 176 "23/23"
Set-RecordedRow 177 "23/23"
Set-RecordedRow 202 "27/30"
Set-RecordedRow 203 "30/30"
Set-RecordedRow 204 "30/30"
Set-RecordedRow 229 "22/26"
Set-RecordedRow 230 "26/26"
Set-RecordedRow 231 "26/26"
Set-RecordedRow 256 "25/28"
Set-RecordedRow 257 "28/28"
Set-RecordedRow 258 "28/28"
Set-RecordedRow 283 "22/26"
Set-RecordedRow 284 "26/26"
Set-RecordedRow 285 "26/26"
Set-RecordedRow 310 "28/29"
Set-RecordedRow 311 "29/29"
Set-RecordedRow 312 "29/29"

# Clean up the scratch cell so it doesn't leak into the saved sheet.
$scratch.Clear()
